$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 2.46
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 3.25
$ws.Range("H4").Value = 2.46
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 3.5
$ws.Range("P4").Value = 1.83
$ws.Range("Q4").Value = 2.02
$ws.Range("AF4").Value = 24
$ws.Range("F5").Value = 4.9
$ws.Range("G5").Value = 5.6
$ws.Range("H5").Value = 1.98
$ws.Range("I5").Value = 2.06
$ws.Range("J5").Value = 3.05
$ws.Range("K5").Value = 3.2
$ws.Range("L5").Value = 1.47
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 2.52
$ws.Range("O5").Value = 1.54
$ws.Range("R5").Value = 1.18
$ws.Range("U5").Value = 1.7
$ws.Range("V5").Value = 1.94
$ws.Range("W5").Value = 1.22
$ws.Range("X5").Value = 8.4
$ws.Range("Y5").Value = 6.8
$ws.Range("Z5").Value = 11
$ws.Range("AA5").Value = 26
$ws.Range("AB5").Value = 13.5
$ws.Range("AE5").Value = 30
$ws.Range("AF5").Value = 38
$ws.Range("AG5").Value = 23
$ws.Range("AH5").Value = 25
$ws.Range("AI5").Value = 70
$ws.Range("AJ5").Value = 160
$ws.Range("AK5").Value = 110
$ws.Range("AL5").Value = 130
$ws.Range("AM5").Value = 370
$ws.Range("AN5").Value = 190
$ws.Range("AO5").Value = 27
$ws.Range("Q6").Value = 2.5
$ws.Range("U6").Value = 1.81
$ws.Range("N8").Value = 1.1
$ws.Range("Q8").Value = 1.86
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 2.3
$ws.Range("I10").Value = 4.8
$ws.Range("V10").Value = 1.26
$ws.Range("H11").Value = 1.69
$ws.Range("S11").Value = 4.3
$ws.Range("U11").Value = 1.72
$ws.Range("J12").Value = 3.4
$ws.Range("N12").Value = 2.74
$ws.Range("P12").Value = 1.89
$ws.Range("S12").Value = 3.25
$ws.Range("U12").Value = 1.05
$ws.Range("F13").Value = 3.55
$ws.Range("H13").Value = 2.2
$ws.Range("I13").Value = 2.24
$ws.Range("J13").Value = 3.55
$ws.Range("P13").Value = 1.86
$ws.Range("R13").Value = 1.33
$ws.Range("S13").Value = 3.6
$ws.Range("T13").Value = 1.79
$ws.Range("U13").Value = 2.08
$ws.Range("Z13").Value = 14.5
$ws.Range("AE13").Value = 24
$ws.Range("F14").Value = 4.5
$ws.Range("H14").Value = 1.94
$ws.Range("I14").Value = 1.96
$ws.Range("J14").Value = 3.6
$ws.Range("P14").Value = 1.87
$ws.Range("S14").Value = 3.55
$ws.Range("T14").Value = 1.83
$ws.Range("U14").Value = 2.04
$ws.Range("V14").Value = 2.04
$ws.Range("X14").Value = 15
$ws.Range("Y14").Value = 8.8
$ws.Range("AH14").Value = 22
$ws.Range("AJ14").Value = 120
$ws.Range("AK14").Value = 50
$ws.Range("AL14").Value = 65
$ws.Range("AN14").Value = 80
$ws.Range("F15").Value = 1.65
$ws.Range("H15").Value = 5.9
$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 4.2
$ws.Range("N15").Value = 4
$ws.Range("O15").Value = 1.28
$ws.Range("R15").Value = 1.41
$ws.Range("S15").Value = 3.1
$ws.Range("V15").Value = 1.18
$ws.Range("AD15").Value = 26
$ws.Range("AG15").Value = 9.6
$ws.Range("AJ15").Value = 19.5
$ws.Range("F16").Value = 1.33
$ws.Range("G16").Value = 1.38
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 5.5
$ws.Range("K16").Value = 6.2
$ws.Range("Q16").Value = 1.62
$ws.Range("S16").Value = 2.52
$ws.Range("U16").Value = 1.86
$ws.Range("V16").Value = 1.09
$ws.Range("W16").Value = 3.6
$ws.Range("X16").Value = 27
$ws.Range("Y16").Value = 1000
$ws.Range("AA16").Value = 460
$ws.Range("AD16").Value = 40
$ws.Range("AH16").Value = 32
$ws.Range("AI16").Value = 150
$ws.Range("AM16").Value = 170
$ws.Range("F17").Value = 2.36
$ws.Range("G17").Value = 3.15
$ws.Range("H17").Value = 2.74
$ws.Range("I17").Value = 4.2
$ws.Range("J17").Value = 2.52
$ws.Range("K17").Value = 4.7
$ws.Range("N17").Value = 1.94
$ws.Range("O17").Value = 1.6
$ws.Range("P17").Value = 1.35
$ws.Range("Q17").Value = 2.66
$ws.Range("S17").Value = 2.66
$ws.Range("V17").Value = 1.32
$ws.Range("W17").Value = 1.46
$ws.Range("X17").Value = 970
$ws.Range("Y17").Value = 970
$ws.Range("AB17").Value = 970
$ws.Range("AC17").Value = 970
$ws.Range("AD17").Value = 970
$ws.Range("AG17").Value = 970
$ws.Range("AH17").Value = 970
$ws.Range("G18").Value = 6.6
$ws.Range("H18").Value = 1.54
$ws.Range("I18").Value = 1.55
$ws.Range("J18").Value = 4.9
$ws.Range("K18").Value = 5.1
$ws.Range("T18").Value = 1.79
$ws.Range("V18").Value = 2.8
$ws.Range("W18").Value = 1.18
$ws.Range("AE18").Value = 19
$ws.Range("AF18").Value = 1000
$ws.Range("L19").Value = 1.46
$ws.Range("M19").Value = 1.01
$ws.Range("N19").Value = 3.15
$ws.Range("O19").Value = 1.38
$ws.Range("P19").Value = 1.79
$ws.Range("R19").Value = 1.25
$ws.Range("S19").Value = 4.2
$ws.Range("T19").Value = 1.01
$ws.Range("U19").Value = 1.01
$ws.Range("V19").Value = 1.35
$ws.Range("W19").Value = 1.72
$ws.Range("X19").Value = 15
$ws.Range("Y19").Value = 15
$ws.Range("Z19").Value = 30
$ws.Range("AA19").Value = 85
$ws.Range("AB19").Value = 10.5
$ws.Range("AC19").Value = 9.2
$ws.Range("AD19").Value = 19
$ws.Range("AE19").Value = 55
$ws.Range("AF19").Value = 17
$ws.Range("AG19").Value = 13.5
$ws.Range("AH19").Value = 23
$ws.Range("AI19").Value = 70
$ws.Range("AJ19").Value = 38
$ws.Range("AK19").Value = 32
$ws.Range("AL19").Value = 50
$ws.Range("AM19").Value = 1000
$ws.Range("AN19").Value = 1000
$ws.Range("AO19").Value = 1000
$ws.Range("F20").Value = 1.43
$ws.Range("K20").Value = 5
$ws.Range("L20").Value = 1.37
$ws.Range("M20").Value = 1.01
$ws.Range("N20").Value = 3.8
$ws.Range("O20").Value = 1.27
$ws.Range("R20").Value = 1.39
$ws.Range("S20").Value = 3.15
$ws.Range("T20").Value = 2.1
$ws.Range("U20").Value = 1.79
$ws.Range("V20").Value = 1.1
$ws.Range("W20").Value = 3.1
$ws.Range("X20").Value = 19.5
$ws.Range("Y20").Value = 30
$ws.Range("Z20").Value = 95
$ws.Range("AA20").Value = 1000
$ws.Range("AB20").Value = 9
$ws.Range("AC20").Value = 12.5
$ws.Range("AD20").Value = 40
$ws.Range("AE20").Value = 1000
$ws.Range("AF20").Value = 9.4
$ws.Range("AG20").Value = 12
$ws.Range("AH20").Value = 32
$ws.Range("AI20").Value = 1000
$ws.Range("AJ20").Value = 14.5
$ws.Range("AK20").Value = 19
$ws.Range("AL20").Value = 48
$ws.Range("AM20").Value = 1000
$ws.Range("AN20").Value = 1000
$ws.Range("AO20").Value = 1000
$ws.Range("G22").Value = 2.3
$ws.Range("H22").Value = 3.7
$ws.Range("I22").Value = 3.95
$ws.Range("G23").Value = 3.25
$ws.Range("H23").Value = 2.6
$ws.Range("I23").Value = 2.72
